$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AutoCompleteSampleSheet")

# Add the new value "Java" to cell A3, right below the existing A2 ("JavaScript")
$ws.Range("A3").Value = "Java"

# Update the selection to reflect the new active cell (A4) as seen in the diff
$ws.Range("A4").Select()
